$wb = $excel.ActiveWorkbook

# --- Analysis Results sheet ---
$ws1 = $wb.Worksheets.Item("Analysis Results")

$ws1.Range("C2").Value = 87.16
$ws1.Range("H2").Value = 90
$ws1.Range("R2").Value = 'Frequency in description matches metadata (monthly)'

$ws1.Range("C3").Value = 76.68
$ws1.Range("H3").Value = 90
$ws1.Range("R3").Value = 'Frequency in description matches metadata (monthly)'

$ws1.Range("Q4").Value = 'No'
$ws1.Range("R4").Value = 'Frequency in description (none) does not match metadata (monthly)'

$ws1.Range("C5").Value = 46.40000000000001
$ws1.Range("E5").Value = 'WHAT, WHY, ESCALATION'
$ws1.Range("H5").Value = 80
$ws1.Range("Q5").Value = 'No'
$ws1.Range("R5").Value = 'Frequency in description (none) does not match metadata (ongoing)'

$ws1.Range("C6").Value = 43.2
$ws1.Range("E6").Value = 'WHAT, WHY, ESCALATION'
$ws1.Range("H6").Value = 80
$ws1.Range("Q6").Value = 'No'
$ws1.Range("R6").Value = 'Frequency in description (none) does not match metadata (daily)'

$ws1.Range("Q7").Value = 'No'
$ws1.Range("R7").Value = 'Frequency in description (none) does not match metadata (monthly)'

$ws1.Range("Q8").Value = 'No'
$ws1.Range("R8").Value = 'Frequency in description (none) does not match metadata (as needed)'

$ws1.Range("Q9").Value = 'No'
$ws1.Range("R9").Value = 'Frequency in description (none) does not match metadata (quarterly)'

$ws1.Range("Q10").Value = 'No'
$ws1.Range("R10").Value = 'Frequency in description (none) does not match metadata (weekly)'

$ws1.Range("Q11").Value = 'No'
$ws1.Range("R11").Value = 'Frequency in description (none) does not match metadata (monthly)'

$ws1.Range("C12").Value = 29.84544
$ws1.Range("H12").Value = 90
$ws1.Range("R12").Value = 'Frequency in description matches metadata (monthly)'

$ws1.Range("C13").Value = 46.68
$ws1.Range("H13").Value = 90
$ws1.Range("R13").Value = 'Frequency in description matches metadata (weekly)'

$ws1.Range("C14").Value = 45.40000000000001
$ws1.Range("H14").Value = 90
$ws1.Range("R14").Value = 'Frequency in description matches metadata (quarterly)'

$ws1.Range("R15").Value = 'Frequency in description (none) does not match metadata (as needed)'

$ws1.Range("C16").Value = 30.22
$ws1.Range("H16").Value = 85
$ws1.Range("Q16").Value = 'No'
$ws1.Range("R16").Value = 'Frequency in description (none) does not match metadata (within 3 days)'

$ws1.Range("Q17").Value = 'No'
$ws1.Range("R17").Value = 'Frequency in description (none) does not match metadata (ongoing)'

$ws1.Range("Q18").Value = 'No'
$ws1.Range("R18").Value = 'Frequency in description (none) does not match metadata (regulatory)'

$ws1.Range("Q19").Value = 'No'
$ws1.Range("R19").Value = 'Frequency in description (none) does not match metadata (ongoing)'

$ws1.Range("Q20").Value = 'No'
$ws1.Range("R20").Value = 'Frequency in description (none) does not match metadata (prompt)'

$ws1.Range("Q21").Value = 'No'
$ws1.Range("R21").Value = 'Frequency in description (none) does not match metadata (internal)'

$ws1.Range("C23").Value = 52.68
$ws1.Range("E23").Value = 'WHEN, WHY'
$ws1.Range("H23").Value = 0

$ws1.Range("C27").Value = 93.18
$ws1.Range("H27").Value = 90
$ws1.Range("R27").Value = 'Frequency in description matches metadata (daily)'

$ws1.Range("C28").Value = 91.8411136
$ws1.Range("H28").Value = 90
$ws1.Range("R28").Value = 'Frequency in description matches metadata (monthly, quarterly)'

$ws1.Range("C29").Value = 89.08
$ws1.Range("H29").Value = 90
$ws1.Range("M29").Value = 'high'
$ws1.Range("O29").Value = 'Yes'
$ws1.Range("R29").Value = 'Frequency in description matches metadata (weekly, ad hoc)'

$ws1.Range("C30").Value = 83.263488
$ws1.Range("H30").Value = 90
$ws1.Range("R30").Value = 'Frequency in description matches metadata (daily)'

$ws1.Range("C31").Value = 95.26173217391305
$ws1.Range("H31").Value = 90
$ws1.Range("R31").Value = 'Frequency in description matches metadata (weekly, monthly)'

$ws1.Columns.Item(8).ColumnWidth = 12.3

# --- Keyword Matches sheet ---
$ws2 = $wb.Worksheets.Item("Keyword Matches")

$ws2.Range("C5").Value = 'by the infosec team'

$ws2.Range("C6").Value = 'by the finance team bef'

$ws2.Range("C12").Value = 'on a monthly basis'

$ws2.Range("C14").Value = 'quarterly'

$ws2.Range("C15").Value = 'as needed'

$ws2.Range("C20").Value = 'None'

$ws2.Range("C23").Value = 'None'

$ws2.Range("C27").Value = 'daily'

$ws2.Range("C29").Value = 'ad hoc, on a weekly basis'

$ws2.Columns.Item(3).ColumnWidth = 28.8

# --- Enhancement Feedback sheet ---
$ws3 = $wb.Worksheets.Item("Enhancement Feedback")

$ws3.Range("C4").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Monthly)'

$ws3.Range("C5").Value = 'Align the frequency in the description with the declared frequency (Ongoing)'

$ws3.Range("C6").Value = 'Align the frequency in the description with the declared frequency (Daily)'

$ws3.Range("C7").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Monthly)'

$ws3.Range("C8").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (As needed)'

$ws3.Range("C9").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Quarterly)'

$ws3.Range("C10").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Weekly)'

$ws3.Range("C11").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Monthly)'

$ws3.Range("C15").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Replace vague timing term ''as needed'' with a specific timeframe or frequency.; Align the frequency in the description with the declared frequency (As needed)'

$ws3.Range("C16").Value = 'Align the frequency in the description with the declared frequency (Within 3 days)'

$ws3.Range("C17").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Ongoing)'

$ws3.Range("C18").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Regulatory)'

$ws3.Range("C19").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Ongoing)'

$ws3.Range("C20").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Prompt)'

$ws3.Range("C21").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).; Align the frequency in the description with the declared frequency (Internal)'

$ws3.Range("C23").Value = 'No specific timing information detected. Add specific frequency (daily, weekly, monthly) or timing (within X days).'

$ws3.Range("C29").Value = 'While ''ad-hoc'' is an allowed frequency, the control would be stronger if it specified what triggers the ad-hoc review.; Multiple frequencies detected. Consider whether this is describing a process rather than a single control.'

# --- Executive Summary sheet ---
$ws5 = $wb.Worksheets.Item("Executive Summary")

$ws5.Range("B4").NumberFormat = "@"
$ws5.Range("B4").Value = '42.3'

$ws5.Range("B15").NumberFormat = "@"
$ws5.Range("B15").Value = '17 (56.7%)'

$ws5.Range("B24").NumberFormat = "@"
$ws5.Range("B24").Value = '15 (50.0%)'

$ws5.Range("B25").NumberFormat = "@"
$ws5.Range("B25").Value = '15 (50.0%)'

Write-Host "Edit applied successfully"